$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1980519480519481
$ws.Range("C2").Value = 0.5714285714285714
$ws.Range("J2").Value = 0.01298701298701299
$ws.Range("P2").Value = 0.1525974025974026
$ws.Range("S2").Value = 0.06493506493506493
# Row 3
$ws.Range("B3").Value = 0.005319148936170213
$ws.Range("C3").Value = 0.02127659574468085
$ws.Range("J3").Value = 0.05319148936170213
$ws.Range("P3").Value = 0.6808510638297872
$ws.Range("S3").Value = 0.2393617021276596
# Row 5
$ws.Range("S5").Value = 1.0
# Row 6
$ws.Range("B6").Value = 0.05529953917050692
$ws.Range("D6").Value = 0.02764976958525346
$ws.Range("F6").Value = 0.05990783410138249
$ws.Range("J6").Value = 0.2903225806451613
$ws.Range("O6").Value = 0.02304147465437788
$ws.Range("Q6").Value = 0.152073732718894
$ws.Range("R6").Value = 0.05990783410138249
$ws.Range("S6").Value = 0.3317972350230415
# Row 7
$ws.Range("B7").Value = 0.1
$ws.Range("D7").Value = 0.005
$ws.Range("F7").Value = 0.055
$ws.Range("J7").Value = 0.16
$ws.Range("O7").Value = 0.015
$ws.Range("Q7").Value = 0.14
$ws.Range("R7").Value = 0.115
$ws.Range("S7").Value = 0.41
# Row 8
$ws.Range("B8").Value = 0.1089108910891089
$ws.Range("D8").Value = 0.007425742574257425
$ws.Range("F8").Value = 0.07425742574257425
$ws.Range("J8").Value = 0.1089108910891089
$ws.Range("O8").Value = 0.03465346534653466
$ws.Range("Q8").Value = 0.1707920792079208
$ws.Range("R8").Value = 0.1188118811881188
$ws.Range("S8").Value = 0.3762376237623762
# Row 9
$ws.Range("B9").Value = 0.1301775147928994
$ws.Range("D9").Value = 0.005917159763313609
$ws.Range("F9").Value = 0.08284023668639054
$ws.Range("J9").Value = 0.136094674556213
$ws.Range("O9").Value = 0.01775147928994083
$ws.Range("Q9").Value = 0.1715976331360947
$ws.Range("R9").Value = 0.0650887573964497
$ws.Range("S9").Value = 0.3905325443786982
# Row 10
$ws.Range("B10").Value = 0.1151419558359622
$ws.Range("D10").Value = 0.01892744479495268
$ws.Range("E10").Value = 0.001577287066246057
$ws.Range("F10").Value = 0.06309148264984227
$ws.Range("J10").Value = 0.1159305993690852
$ws.Range("O10").Value = 0.01971608832807571
$ws.Range("Q10").Value = 0.2208201892744479
$ws.Range("R10").Value = 0.07886435331230283
$ws.Range("S10").Value = 0.3659305993690852
# Row 11
$ws.Range("G11").Value = 0.1442006269592477
$ws.Range("J11").Value = 0.1128526645768025
$ws.Range("K11").Value = 0.2068965517241379
$ws.Range("L11").Value = 0.5297805642633229
$ws.Range("S11").Value = 0.006269592476489028
# Row 12
$ws.Range("G12").Value = 0.7588235294117647
$ws.Range("J12").Value = 0.2
$ws.Range("L12").Value = 0.01176470588235294
$ws.Range("S12").Value = 0.02941176470588235
# Row 13
$ws.Range("G13").Value = 0.6170212765957447
$ws.Range("J13").Value = 0.2978723404255319
$ws.Range("S13").Value = 0.0851063829787234
# Row 15
$ws.Range("F15").Value = 0.01951219512195122
$ws.Range("H15").Value = 0.175609756097561
$ws.Range("I15").Value = 0.07804878048780488
$ws.Range("J15").Value = 0.2829268292682927
$ws.Range("K15").Value = 0.07804878048780488
$ws.Range("M15").Value = 0.02439024390243903
$ws.Range("O15").Value = 0.03902439024390244
$ws.Range("S15").Value = 0.3024390243902439
# Row 16
$ws.Range("F16").Value = 0.03684210526315789
$ws.Range("H16").Value = 0.1157894736842105
$ws.Range("I16").Value = 0.05789473684210526
$ws.Range("J16").Value = 0.4315789473684211
$ws.Range("K16").Value = 0.1157894736842105
$ws.Range("M16").Value = 0.03157894736842105
$ws.Range("O16").Value = 0.04736842105263158
$ws.Range("S16").Value = 0.1631578947368421
# Row 17
$ws.Range("F17").Value = 0.02528735632183908
$ws.Range("H17").Value = 0.1655172413793103
$ws.Range("I17").Value = 0.07816091954022988
$ws.Range("J17").Value = 0.4413793103448276
$ws.Range("K17").Value = 0.09195402298850575
$ws.Range("M17").Value = 0.01149425287356322
$ws.Range("O17").Value = 0.04367816091954023
$ws.Range("S17").Value = 0.1425287356321839
# Row 18
$ws.Range("F18").Value = 0.01538461538461539
$ws.Range("H18").Value = 0.1692307692307692
$ws.Range("I18").Value = 0.08205128205128205
$ws.Range("J18").Value = 0.4051282051282051
$ws.Range("K18").Value = 0.08205128205128205
$ws.Range("M18").Value = 0.02051282051282051
$ws.Range("O18").Value = 0.06153846153846154
$ws.Range("S18").Value = 0.1641025641025641
# Row 19
$ws.Range("F19").Value = 0.01813880126182965
$ws.Range("H19").Value = 0.1924290220820189
$ws.Range("I19").Value = 0.07334384858044164
$ws.Range("J19").Value = 0.3698738170347003
$ws.Range("K19").Value = 0.1214511041009464
$ws.Range("M19").Value = 0.02287066246056782
$ws.Range("N19").Value = 0.001577287066246057
$ws.Range("O19").Value = 0.06309148264984227
$ws.Range("S19").Value = 0.1372239747634069

$wb.Save()
